{"js": "// Replace the \"A general \" run text with the new sentence start,\n// and collapse the \"IEEE style guide is available at \" run to a single space.\n\nconst body = context.document.body;\n\nconst firstSearch = body.search(\"A general \", { matchCase: true, matchWholeWord: false });\nfirstSearch.load(\"items\");\nawait context.sync();\n\nif (firstSearch.items.length > 0) {\n  firstSearch.items[0].insertText(\"T twin notch, notch, band pass, \", \"Replace\");\n  await context.sync();\n}\n\nconst secondSearch = body.search(\"IEEE style guide is available at \", { matchCase: true, matchWholeWord: false });\nsecondSearch.load(\"items\");\nawait context.sync();\n\nif (secondSearch.items.length > 0) {\n  secondSearch.items[0].insertText(\" \", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the first run's text \"A general \" with the new lead-in text.\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Execute(\"A general \", $false, $false, $false, $false, $false, $true, 1, $false, \"T twin notch, notch, band pass, \", 2)\n\n# Collapse the second run's text \"IEEE style guide is available at \" down to a single space,\n# keeping that run's own formatting (color) intact.\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"IEEE style guide is available at \", $false, $false, $false, $false, $false, $true, 1, $false, \" \", 2)\n"}
